$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38: B38 changes from text "3" to numeric 3
$ws.Range("B38").Value = 3

# New row 39 appended
$ws.Range("A39").Value = "Ruilin"

# B39 must remain a text string "3" (not auto-converted to a number)
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "3"
$ws.Range("B39").ClearFormats()

$ws.Range("C39").Value = "无"
$ws.Range("D39").Value = "DIS"
$ws.Range("E39").Value = "RES"
$ws.Range("F39").Value = "df7b0ece-3727-4ec6-95ce-2a2839e398ed"
$ws.Range("G39").Value = "SkhQHMW0W_annotated.xlsx"
$ws.Range("H39").Value = "This is necessary to get an idea of total amount of communication that was sufficient to reach perplexity 72.24 at the end of 40-th epoch."
